# Applies the Q3 "Worst" output updates (per commit: "Updated by VB & SB / Question 3 completed")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Cells = @{ "D"=71.40000000000001; "E"=1998; "G"=1998; "H"=356785.71; "J"=178.57; "K"=0.4 } },
    @{ Row = 3; Cells = @{ "D"=53.8; "E"=1998; "F"=599; "G"=1398; "H"=376858.46; "I"=29.98; "J"=188.62; "K"=0.4 } },
    @{ Row = 4; Cells = @{ "D"=71.40000000000001; "E"=1998; "G"=1998; "H"=356785.71; "J"=178.57; "K"=0.4 } },
    @{ Row = 5; Cells = @{ "D"=66.59999999999999; "E"=1998; "G"=1998; "H"=333000; "J"=166.67; "K"=0.3 } },
    @{ Row = 6; Cells = @{ "D"=71.40000000000001; "E"=1998; "G"=1998; "H"=356785.71; "J"=178.57; "K"=0.4 } },
    @{ Row = 7; Cells = @{ "D"=53.8; "E"=1998; "F"=599; "G"=1398; "H"=376858.46; "I"=29.98; "J"=188.62 } },
    @{ Row = 8; Cells = @{ "D"=58.3; "E"=1998; "F"=599; "G"=1399; "H"=399278.33; "I"=29.98; "J"=199.84; "K"=0.4 } },
    @{ Row = 9; Cells = @{ "D"=53.8; "E"=1998; "F"=599; "G"=1398; "H"=376858.46; "I"=29.98; "J"=188.62; "K"=0.4 } },
    @{ Row = 10; Cells = @{ "D"=71.40000000000001; "E"=1998; "G"=1998; "H"=356785.71; "J"=178.57; "K"=0.4 } },
    @{ Row = 11; Cells = @{ "D"=53.8; "E"=1998; "F"=599; "G"=1398; "H"=376858.46; "I"=29.98; "J"=188.62; "K"=0.4 } },
    @{ Row = 12; Cells = @{ "D"=58.3; "E"=1998; "F"=599; "G"=1399; "H"=399278.33; "I"=29.98; "J"=199.84; "K"=0.4 } },
    @{ Row = 13; Cells = @{ "D"=58.3; "E"=1998; "F"=599; "G"=1399; "H"=399278.33; "J"=199.84; "K"=0.4 } },
    @{ Row = 14; Cells = @{ "D"=42.8; "E"=1998; "F"=799; "G"=1199; "H"=316114.88; "I"=39.99; "J"=158.22; "K"=0.3 } },
    @{ Row = 15; Cells = @{ "D"=46.1; "E"=1998; "F"=799; "G"=1198; "H"=331212.18; "J"=165.77; "K"=0.3 } },
    @{ Row = 16; Cells = @{ "D"=42.8; "E"=1998; "F"=799; "G"=1199; "H"=316114.88; "I"=39.99; "J"=158.22; "K"=0.3 } },
    @{ Row = 17; Cells = @{ "D"=40; "E"=1998; "F"=799; "G"=1199; "H"=303030.56; "I"=39.99; "J"=151.67 } },
    @{ Row = 18; Cells = @{ "D"=42.8; "E"=1998; "F"=799; "G"=1199; "H"=316114.88; "I"=39.99; "J"=158.22; "K"=0.3 } },
    @{ Row = 19; Cells = @{ "D"=46.1; "E"=1998; "F"=799; "G"=1198; "H"=331212.18; "I"=39.99; "J"=165.77; "K"=0.3 } },
    @{ Row = 20; Cells = @{ "D"=50; "E"=1998; "F"=799; "G"=1198; "H"=348825.69; "J"=174.59; "K"=0.3 } },
    @{ Row = 21; Cells = @{ "D"=46.1; "E"=1998; "F"=799; "G"=1198; "H"=331212.18; "I"=39.99; "J"=165.77; "K"=0.3 } },
    @{ Row = 22; Cells = @{ "D"=42.8; "E"=1998; "F"=799; "G"=1199; "H"=316114.88; "I"=39.99; "J"=158.22; "K"=0.3 } },
    @{ Row = 23; Cells = @{ "D"=46.1; "E"=1998; "F"=799; "G"=1198; "H"=331212.18; "I"=39.99; "J"=165.77; "K"=0.3 } },
    @{ Row = 24; Cells = @{ "D"=50; "E"=1998; "F"=799; "G"=1198; "H"=348825.69; "I"=39.99; "J"=174.59; "K"=0.3 } },
    @{ Row = 25; Cells = @{ "D"=50; "E"=1998; "F"=799; "G"=1198; "H"=348825.69; "I"=39.99; "J"=174.59; "K"=0.3 } },
    @{ Row = 26; Cells = @{ "D"=71.40000000000001; "E"=1998; "G"=1998; "H"=315160.71; "J"=157.74; "K"=0.3 } },
    @{ Row = 27; Cells = @{ "D"=0; "E"=1998; "F"=1998; "G"=0; "H"=319680; "I"=100; "J"=160 } },
    @{ Row = 28; Cells = @{ "D"=71.40000000000001; "E"=1998; "G"=1998; "H"=315160.71; "J"=157.74 } },
    @{ Row = 29; Cells = @{ "D"=66.59999999999999; "E"=1998; "G"=1998; "H"=294150; "J"=147.22 } },
    @{ Row = 30; Cells = @{ "D"=71.40000000000001; "E"=1998; "G"=1998; "H"=315160.71; "J"=157.74 } },
    @{ Row = 31; Cells = @{ "D"=0; "E"=1998; "F"=1998; "G"=0; "H"=319680; "I"=100; "J"=160; "K"=0.3 } },
    @{ Row = 32; Cells = @{ "E"=1998; "F"=1998; "H"=319680; "K"=0.3 } },
    @{ Row = 33; Cells = @{ "D"=0; "E"=1998; "F"=1998; "G"=0; "H"=319680; "I"=100; "J"=160; "K"=0.3 } },
    @{ Row = 34; Cells = @{ "D"=71.40000000000001; "E"=1998; "G"=1998; "H"=315160.71; "J"=157.74; "K"=0.3 } },
    @{ Row = 35; Cells = @{ "D"=0; "E"=1998; "F"=1998; "G"=0; "H"=319680; "I"=100; "J"=160; "K"=0.3 } },
    @{ Row = 36; Cells = @{ "E"=1998; "F"=1998; "H"=319680; "K"=0.3 } }
)

foreach ($update in $updates) {
    $r = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $ws.Range("$col$r").Value = $update.Cells[$col]
    }
}